# The "#测试说明" (test-notes) column D is being removed from the hero
# config sheet. Deleting the whole column shifts columns E:H left into
# D:G, carrying their values/number formats/styles with them (this is
# exactly what the OOXML diff shows: the old E/F/G/H content now lives in
# D/E/F/G, the orphaned "#测试说明" shared string disappears, and the
# sheet's dimension shrinks from A1:H8 to A1:G8).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D").Delete()

# Match the author's final selection/cursor position recorded in the
# worksheet view.
$ws.Range("E14").Select() | Out-Null
